$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.233.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.609.33'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.00'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.87'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.636.16'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.55'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.89%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.073.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.113.05'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.45'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.619.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.67%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.94'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.94'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.521'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +10.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.36'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.97'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +8.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0805'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.35%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.38'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '162.50'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.58'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +16.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.29'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.30%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.62'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.95'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.79%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.854'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '304.32'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '134.46'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.03'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.05%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0551'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0243'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.07%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.99'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +10.75%  '
